$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix department code value: drop the stray trailing space ("D001 " -> "D001")
$ws.Range("A2").Value = "D001"

# 2. Remove the mailto hyperlink on D2 (keep the display text as plain text)
$ws.Range("D2").Hyperlinks.Delete()

# 3. Recolor the body of the table from the light-blue fill to white
#    (xlThemeColorLight1 = 2 -> theme="0"/Background 1, matches the target fill)
#    Only the cells that actually carried the light-blue fill (plus D4, which
#    picks up the same fill the hyperlink-style column already used) change;
#    B2/C2/A4/B4/C4 keep their "no fill" formatting untouched.
$ws.Range("A2").Interior.ThemeColor = 2
$ws.Range("D2").Interior.ThemeColor = 2
$ws.Range("A3:D3").Interior.ThemeColor = 2
$ws.Range("D4").Interior.ThemeColor = 2

# 4. Update the active selection to B3
$ws.Range("B3").Select()
